$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells for the season record columns, matching the
# existing bold/centered/bordered header style used by A1:AC1.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$newHeaders = $ws.Range("AD1:AF1")
$newHeaders.Font.Bold = $true
$newHeaders.HorizontalAlignment = -4108
$newHeaders.VerticalAlignment = -4160
$newHeaders.Borders.LineStyle = 1

# Fill in the season record (Wins/Losses/Ties) for every player row.
$wins = 103
$losses = 58
$ties = 1

for ($row = 2; $row -le 47; $row++) {
    $ws.Cells.Item($row, 30).Value = $wins    # AD
    $ws.Cells.Item($row, 31).Value = $losses  # AE
    $ws.Cells.Item($row, 32).Value = $ties    # AF
}
